$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K5").Select() | Out-Null
$ws.Range("K5").ClearContents() | Out-Null
